# Scheduled runner update: refresh computed price/profit figures on each
# class sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect latest
# market data pulls.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1469.2727
$ws.Range("J129").Value = 5609
$ws.Range("L129").Value = 16827
$ws.Range("N129").Value = -26827
$ws.Range("H137").Value = 20834452
$ws.Range("I137").Value = 3547083.2
$ws.Range("K137").Value = 10641249.6
$ws.Range("M137").Value = -10638699.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1226855.5
$ws.Range("I2").Value = 1265.2
$ws.Range("J2").Value = 2102277
$ws.Range("K2").Value = 1265.2
$ws.Range("L2").Value = 2102277
$ws.Range("M2").Value = -1152.2
$ws.Range("N2").Value = -2102503
$ws.Range("H32").Value = 5511.453
$ws.Range("I32").Value = 5511.8228
$ws.Range("K32").Value = 5511.8228
$ws.Range("M32").Value = -5224.8228
$ws.Range("H39").Value = 1650
$ws.Range("I39").Value = 1500
$ws.Range("J39").Value = 1800
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 1800
$ws.Range("M39").Value = -980
$ws.Range("N39").Value = -2840
$ws.Range("H42").Value = 12600
$ws.Range("J42").Value = 12600
$ws.Range("L42").Value = 12600
$ws.Range("N42").Value = -13572
$ws.Range("H61").Value = 10419404
$ws.Range("I61").Value = 12348532
$ws.Range("K61").Value = 12348532
$ws.Range("M61").Value = -12348320
$ws.Range("H74").Value = 2987700.2
$ws.Range("I74").Value = 3846618
$ws.Range("J74").Value = 10118.333
$ws.Range("K74").Value = 3846618
$ws.Range("L74").Value = 10118.333
$ws.Range("M74").Value = -3845744
$ws.Range("N74").Value = -11866.333
$ws.Range("H77").Value = 2987700.2
$ws.Range("I77").Value = 3846618
$ws.Range("J77").Value = 10118.333
$ws.Range("K77").Value = 19233090
$ws.Range("L77").Value = 50591.665
$ws.Range("M77").Value = -19228722
$ws.Range("N77").Value = -59327.665
$ws.Range("H116").Value = 1226855.5
$ws.Range("I116").Value = 1265.2
$ws.Range("J116").Value = 2102277
$ws.Range("K116").Value = 1265.2
$ws.Range("L116").Value = 2102277
$ws.Range("M116").Value = 1028.8
$ws.Range("N116").Value = -2106865
$ws.Range("H132").Value = 867305.7
$ws.Range("I132").Value = 1151361.8
$ws.Range("K132").Value = 3454085.4
$ws.Range("M132").Value = -3451555.4
$ws.Range("H136").Value = 10419404
$ws.Range("I136").Value = 12348532
$ws.Range("K136").Value = 37045596
$ws.Range("M136").Value = -37043046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1226855.5
$ws.Range("I3").Value = 1265.2
$ws.Range("J3").Value = 2102277
$ws.Range("K3").Value = 1265.2
$ws.Range("L3").Value = 2102277
$ws.Range("M3").Value = -1151.2
$ws.Range("N3").Value = -2102505
$ws.Range("H38").Value = 17000
$ws.Range("J38").Value = 17000
$ws.Range("L38").Value = 17000
$ws.Range("N38").Value = -17832
$ws.Range("H134").Value = 61107.85
$ws.Range("I134").Value = 76071.81
$ws.Range("J134").Value = 1252
$ws.Range("K134").Value = 228215.43
$ws.Range("L134").Value = 3756
$ws.Range("M134").Value = -225680.43
$ws.Range("N134").Value = -8826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2655
$ws.Range("I16").Value = 2668.5
$ws.Range("J16").Value = 2622.6
$ws.Range("K16").Value = 2668.5
$ws.Range("L16").Value = 2622.6
$ws.Range("M16").Value = -2381.5
$ws.Range("N16").Value = -3196.6
$ws.Range("H31").Value = 2272.4814
$ws.Range("I31").Value = 1328.5834
$ws.Range("J31").Value = 3027.6
$ws.Range("K31").Value = 1328.5834
$ws.Range("L31").Value = 3027.6
$ws.Range("M31").Value = -1033.5834
$ws.Range("N31").Value = -3617.6
$ws.Range("H34").Value = 2272.4814
$ws.Range("I34").Value = 1328.5834
$ws.Range("J34").Value = 3027.6
$ws.Range("K34").Value = 1328.5834
$ws.Range("L34").Value = 3027.6
$ws.Range("M34").Value = -1126.5834
$ws.Range("N34").Value = -3431.6
$ws.Range("H113").Value = 2655
$ws.Range("I113").Value = 2668.5
$ws.Range("J113").Value = 2622.6
$ws.Range("K113").Value = 2668.5
$ws.Range("L113").Value = 2622.6
$ws.Range("M113").Value = -498.5
$ws.Range("N113").Value = -6962.6
$ws.Range("H122").Value = 1151.3529
$ws.Range("I122").Value = 1005.3333
$ws.Range("J122").Value = 1315.625
$ws.Range("K122").Value = 3015.9999
$ws.Range("L122").Value = 3946.875
$ws.Range("M122").Value = -565.9998999999998
$ws.Range("N122").Value = -8846.875
$ws.Range("H132").Value = 2212.8167
$ws.Range("I132").Value = 2154.843
$ws.Range("J132").Value = 2541.3333
$ws.Range("K132").Value = 6464.529
$ws.Range("L132").Value = 7623.999899999999
$ws.Range("M132").Value = -3934.529
$ws.Range("N132").Value = -12683.9999
$ws.Range("H134").Value = 1678.8667
$ws.Range("I134").Value = 1745.3334
$ws.Range("K134").Value = 5236.0002
$ws.Range("M134").Value = -2701.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 553
$ws.Range("I113").Value = 532
$ws.Range("J113").Value = 588
$ws.Range("K113").Value = 1596
$ws.Range("L113").Value = 1764
$ws.Range("M113").Value = 574
$ws.Range("N113").Value = -6104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1601.6
$ws.Range("I102").Value = 1493.375
$ws.Range("J102").Value = 1725.2858
$ws.Range("K102").Value = 1493.375
$ws.Range("L102").Value = 1725.2858
$ws.Range("M102").Value = 128.625
$ws.Range("N102").Value = -4969.2858
$ws.Range("H122").Value = 3945.7437
$ws.Range("I122").Value = 3032.375
$ws.Range("J122").Value = 5407.1333
$ws.Range("K122").Value = 9097.125
$ws.Range("L122").Value = 16221.3999
$ws.Range("M122").Value = -6647.125
$ws.Range("N122").Value = -21121.3999
$ws.Range("H132").Value = 1732.491
$ws.Range("I132").Value = 1603.9688
$ws.Range("J132").Value = 1911.3043
$ws.Range("K132").Value = 4811.9064
$ws.Range("L132").Value = 5733.9129
$ws.Range("M132").Value = -2281.9064
$ws.Range("N132").Value = -10793.9129

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1318.3529
$ws.Range("I40").Value = 1537.4546
$ws.Range("J40").Value = 916.6667
$ws.Range("K40").Value = 1537.4546
$ws.Range("L40").Value = 916.6667
$ws.Range("M40").Value = -1401.4546
$ws.Range("N40").Value = -1188.6667
$ws.Range("H61").Value = 2004.4546
$ws.Range("I61").Value = 2143.625
$ws.Range("J61").Value = 1633.3334
$ws.Range("K61").Value = 2143.625
$ws.Range("L61").Value = 1633.3334
$ws.Range("M61").Value = -1941.625
$ws.Range("N61").Value = -2037.3334
$ws.Range("H113").Value = 2004.4546
$ws.Range("I113").Value = 2143.625
$ws.Range("J113").Value = 1633.3334
$ws.Range("K113").Value = 2143.625
$ws.Range("L113").Value = 1633.3334
$ws.Range("M113").Value = 26.375
$ws.Range("N113").Value = -5973.3334
$ws.Range("H122").Value = 1608.3334
$ws.Range("I122").Value = 1630
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4890
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -2440
$ws.Range("N122").Value = -9400
$ws.Range("H136").Value = 4715.357
$ws.Range("I136").Value = 4721.2
$ws.Range("K136").Value = 14163.6
$ws.Range("M136").Value = -11613.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2227.7666
$ws.Range("I122").Value = 1399.0555
$ws.Range("J122").Value = 3470.8333
$ws.Range("K122").Value = 4197.166499999999
$ws.Range("L122").Value = 10412.4999
$ws.Range("M122").Value = -1747.166499999999
$ws.Range("N122").Value = -15312.4999
$ws.Range("H136").Value = 10765.9
$ws.Range("I136").Value = 11795.444
$ws.Range("K136").Value = 35386.33199999999
$ws.Range("M136").Value = -32836.33199999999
